$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.37596266666667
$ws.Range("H2").Value = 34.127888
$ws.Range("I2").Value = 0.05604480707695051
$ws.Range("J2").Value = 0.05604480707695052
$ws.Range("M2").Value = 19.21315233333334
$ws.Range("N2").Value = 57.63945700000001
$ws.Range("O2").Value = 0.04451179209991234
$ws.Range("P2").Value = 0.04451179209991233
$ws.Range("Q2").Value = 218.5681036529796
$ws.Range("R2").Value = 1967.112932876816
$ws.Range("S2").Value = 0.002494654800888917
$ws.Range("T2").Value = 0.002494654800888917
$ws.Range("G3").Value = 11.37596266666667
$ws.Range("H3").Value = 34.127888
$ws.Range("I3").Value = 0.05604480707695051
$ws.Range("J3").Value = 0.05604480707695052
$ws.Range("O3").Value = 0.2141755495962477
$ws.Range("P3").Value = 0.2141755495962477
$ws.Range("Q3").Value = 1051.675107104455
$ws.Range("R3").Value = 9465.075963940095
$ws.Range("S3").Value = 0.01200342735772155
$ws.Range("T3").Value = 0.01200342735772155
$ws.Range("G4").Value = 11.37596266666667
$ws.Range("H4").Value = 34.127888
$ws.Range("I4").Value = 0.05604480707695051
$ws.Range("J4").Value = 0.05604480707695052
$ws.Range("M4").Value = 166.8580016666666
$ws.Range("N4").Value = 500.5740049999999
$ws.Range("O4").Value = 0.3865658561145097
$ws.Range("P4").Value = 0.3865658561145097
$ws.Range("Q4").Value = 1898.170397594604
$ws.Range("R4").Value = 17083.53357835144
$ws.Range("S4").Value = 0.02166500882847391
$ws.Range("T4").Value = 0.02166500882847391
$ws.Range("G5").Value = 11.37596266666667
$ws.Range("H5").Value = 34.127888
$ws.Range("I5").Value = 0.05604480707695051
$ws.Range("J5").Value = 0.05604480707695052
$ws.Range("M5").Value = 41.09915599999999
$ws.Range("N5").Value = 123.297468
$ws.Range("O5").Value = 0.09521587377309249
$ws.Range("P5").Value = 0.09521587377309249
$ws.Range("Q5").Value = 467.5424642875092
$ws.Range("R5").Value = 4207.882178587583
$ws.Range("S5").Value = 0.00533635527627624
$ws.Range("T5").Value = 0.005336355276276241
$ws.Range("G6").Value = 11.37596266666667
$ws.Range("H6").Value = 34.127888
$ws.Range("I6").Value = 0.05604480707695051
$ws.Range("J6").Value = 0.05604480707695052
$ws.Range("M6").Value = 112.0244103333333
$ws.Range("N6").Value = 336.073231
$ws.Range("O6").Value = 0.2595309284162377
$ws.Range("P6").Value = 0.2595309284162377
$ws.Range("Q6").Value = 1274.385509707348
$ws.Range("R6").Value = 11469.46958736613
$ws.Range("S6").Value = 0.0145453608135899
$ws.Range("T6").Value = 0.0145453608135899
$ws.Range("I7").Value = 0.765548861900355
$ws.Range("J7").Value = 0.7655488619003551
$ws.Range("M7").Value = 19.21315233333334
$ws.Range("N7").Value = 57.63945700000001
$ws.Range("O7").Value = 0.04451179209991234
$ws.Range("P7").Value = 0.04451179209991233
$ws.Range("Q7").Value = 2985.549807844245
$ws.Range("R7").Value = 26869.9482705982
$ws.Range("S7").Value = 0.0340759517832331
$ws.Range("T7").Value = 0.0340759517832331
$ws.Range("I8").Value = 0.765548861900355
$ws.Range("J8").Value = 0.7655488619003551
$ws.Range("O8").Value = 0.2141755495962477
$ws.Range("P8").Value = 0.2141755495962477
$ws.Range("S8").Value = 0.1639618482402905
$ws.Range("T8").Value = 0.1639618482402905
$ws.Range("I9").Value = 0.765548861900355
$ws.Range("J9").Value = 0.7655488619003551
$ws.Range("M9").Value = 166.8580016666666
$ws.Range("N9").Value = 500.5740049999999
$ws.Range("O9").Value = 0.3865658561145097
$ws.Range("P9").Value = 0.3865658561145097
$ws.Range("Q9").Value = 25928.22178112423
$ws.Range("R9").Value = 233353.9960301181
$ws.Range("S9").Value = 0.2959350511979993
$ws.Range("T9").Value = 0.2959350511979993
$ws.Range("I10").Value = 0.765548861900355
$ws.Range("J10").Value = 0.7655488619003551
$ws.Range("M10").Value = 41.09915599999999
$ws.Range("N10").Value = 123.297468
$ws.Range("O10").Value = 0.09521587377309249
$ws.Range("P10").Value = 0.09521587377309249
$ws.Range("Q10").Value = 6386.436497746357
$ws.Range("R10").Value = 57477.92847971722
$ws.Range("S10").Value = 0.07289240380183881
$ws.Range("T10").Value = 0.07289240380183883
$ws.Range("I11").Value = 0.765548861900355
$ws.Range("J11").Value = 0.7655488619003551
$ws.Range("M11").Value = 112.0244103333333
$ws.Range("N11").Value = 336.073231
$ws.Range("O11").Value = 0.2595309284162377
$ws.Range("P11").Value = 0.2595309284162377
$ws.Range("Q11").Value = 17407.57846198385
$ws.Range("R11").Value = 156668.2061578547
$ws.Range("S11").Value = 0.1986836068769933
$ws.Range("T11").Value = 0.1986836068769933
$ws.Range("G12").Value = 11.89345866666667
$ws.Range("H12").Value = 35.680376
$ws.Range("I12").Value = 0.05859430238850571
$ws.Range("J12").Value = 0.05859430238850571
$ws.Range("M12").Value = 19.21315233333334
$ws.Range("N12").Value = 57.63945700000001
$ws.Range("O12").Value = 0.04451179209991234
$ws.Range("P12").Value = 0.04451179209991233
$ws.Range("Q12").Value = 228.5108331328703
$ws.Range("R12").Value = 2056.597498195832
$ws.Range("S12").Value = 0.002608137406156563
$ws.Range("T12").Value = 0.002608137406156563
$ws.Range("G13").Value = 11.89345866666667
$ws.Range("H13").Value = 35.680376
$ws.Range("I13").Value = 0.05859430238850571
$ws.Range("J13").Value = 0.05859430238850571
$ws.Range("O13").Value = 0.2141755495962477
$ws.Range("P13").Value = 0.2141755495962477
$ws.Range("Q13").Value = 1099.516127435933
$ws.Range("R13").Value = 9895.645146923393
$ws.Range("S13").Value = 0.01254946691726694
$ws.Range("T13").Value = 0.01254946691726694
$ws.Range("G14").Value = 11.89345866666667
$ws.Range("H14").Value = 35.680376
$ws.Range("I14").Value = 0.05859430238850571
$ws.Range("J14").Value = 0.05859430238850571
$ws.Range("M14").Value = 166.8580016666666
$ws.Range("N14").Value = 500.5740049999999
$ws.Range("O14").Value = 0.3865658561145097
$ws.Range("P14").Value = 0.3865658561145097
$ws.Range("Q14").Value = 1984.518746025097
$ws.Range("R14").Value = 17860.66871422588
$ws.Range("S14").Value = 0.02265055666624517
$ws.Range("T14").Value = 0.02265055666624517
$ws.Range("G15").Value = 11.89345866666667
$ws.Range("H15").Value = 35.680376
$ws.Range("I15").Value = 0.05859430238850571
$ws.Range("J15").Value = 0.05859430238850571
$ws.Range("M15").Value = 41.09915599999999
$ws.Range("N15").Value = 123.297468
$ws.Range("O15").Value = 0.09521587377309249
$ws.Range("P15").Value = 0.09521587377309249
$ws.Range("Q15").Value = 488.8111131208853
$ws.Range("R15").Value = 4399.300018087968
$ws.Range("S15").Value = 0.005579107700046371
$ws.Range("T15").Value = 0.005579107700046371
$ws.Range("G16").Value = 11.89345866666667
$ws.Range("H16").Value = 35.680376
$ws.Range("I16").Value = 0.05859430238850571
$ws.Range("J16").Value = 0.05859430238850571
$ws.Range("M16").Value = 112.0244103333333
$ws.Range("N16").Value = 336.073231
$ws.Range("O16").Value = 0.2595309284162377
$ws.Range("P16").Value = 0.2595309284162377
$ws.Range("Q16").Value = 1332.357693957206
$ws.Range("R16").Value = 11991.21924561486
$ws.Range("S16").Value = 0.01520703369879066
$ws.Range("T16").Value = 0.01520703369879066
$ws.Range("G17").Value = 23.69116533333333
$ws.Range("H17").Value = 71.07349600000001
$ws.Range("I17").Value = 0.1167168730630039
$ws.Range("J17").Value = 0.1167168730630039
$ws.Range("M17").Value = 19.21315233333334
$ws.Range("N17").Value = 57.63945700000001
$ws.Range("O17").Value = 0.04451179209991234
$ws.Range("P17").Value = 0.04451179209991233
$ws.Range("Q17").Value = 455.1819685035192
$ws.Range("R17").Value = 4096.637716531673
$ws.Range("S17").Value = 0.005195277188332288
$ws.Range("T17").Value = 0.005195277188332288
$ws.Range("G18").Value = 23.69116533333333
$ws.Range("H18").Value = 71.07349600000001
$ws.Range("I18").Value = 0.1167168730630039
$ws.Range("J18").Value = 0.1167168730630039
$ws.Range("O18").Value = 0.2141755495962477
$ws.Range("P18").Value = 0.2141755495962477
$ws.Range("Q18").Value = 2190.18025721627
$ws.Range("R18").Value = 19711.62231494643
$ws.Range("S18").Value = 0.02499790043542434
$ws.Range("T18").Value = 0.02499790043542435
$ws.Range("G19").Value = 23.69116533333333
$ws.Range("H19").Value = 71.07349600000001
$ws.Range("I19").Value = 0.1167168730630039
$ws.Range("J19").Value = 0.1167168730630039
$ws.Range("M19").Value = 166.8580016666666
$ws.Range("N19").Value = 500.5740049999999
$ws.Range("O19").Value = 0.3865658561145097
$ws.Range("P19").Value = 0.3865658561145097
$ws.Range("Q19").Value = 3953.060504674608
$ws.Range("R19").Value = 35577.54454207148
$ws.Range("S19").Value = 0.04511875795860866
$ws.Range("T19").Value = 0.04511875795860867
$ws.Range("G20").Value = 23.69116533333333
$ws.Range("H20").Value = 71.07349600000001
$ws.Range("I20").Value = 0.1167168730630039
$ws.Range("J20").Value = 0.1167168730630039
$ws.Range("M20").Value = 41.09915599999999
$ws.Range("N20").Value = 123.297468
$ws.Range("O20").Value = 0.09521587377309249
$ws.Range("P20").Value = 0.09521587377309249
$ws.Range("Q20").Value = 973.6868998564586
$ws.Range("R20").Value = 8763.182098708128
$ws.Range("S20").Value = 0.01111329905275704
$ws.Range("T20").Value = 0.01111329905275704
$ws.Range("G21").Value = 23.69116533333333
$ws.Range("H21").Value = 71.07349600000001
$ws.Range("I21").Value = 0.1167168730630039
$ws.Range("J21").Value = 0.1167168730630039
$ws.Range("M21").Value = 112.0244103333333
$ws.Range("N21").Value = 336.073231
$ws.Range("O21").Value = 0.2595309284162377
$ws.Range("P21").Value = 0.2595309284162377
$ws.Range("Q21").Value = 2653.988826576175
$ws.Range("R21").Value = 23885.89943918558
$ws.Range("S21").Value = 0.03029163842788157
$ws.Range("T21").Value = 0.03029163842788157
$ws.Range("G22").Value = 0.628254
$ws.Range("H22").Value = 1.884762
$ws.Range("I22").Value = 0.003095155571184698
$ws.Range("J22").Value = 0.003095155571184698
$ws.Range("M22").Value = 19.21315233333334
$ws.Range("N22").Value = 57.63945700000001
$ws.Range("O22").Value = 0.04451179209991234
$ws.Range("P22").Value = 0.04451179209991233
$ws.Range("Q22").Value = 12.070739806026
$ws.Range("R22").Value = 108.636658254234
$ws.Range("S22").Value = 0.0001377709213014587
$ws.Range("T22").Value = 0.0001377709213014587
$ws.Range("G23").Value = 0.628254
$ws.Range("H23").Value = 1.884762
$ws.Range("I23").Value = 0.003095155571184698
$ws.Range("J23").Value = 0.003095155571184698
$ws.Range("O23").Value = 0.2141755495962477
$ws.Range("P23").Value = 0.2141755495962477
$ws.Range("Q23").Value = 58.080279629856
$ws.Range("R23").Value = 522.722516668704
$ws.Range("S23").Value = 0.0006629066455443707
$ws.Range("T23").Value = 0.0006629066455443707
$ws.Range("G24").Value = 0.628254
$ws.Range("H24").Value = 1.884762
$ws.Range("I24").Value = 0.003095155571184698
$ws.Range("J24").Value = 0.003095155571184698
$ws.Range("M24").Value = 166.8580016666666
$ws.Range("N24").Value = 500.5740049999999
$ws.Range("O24").Value = 0.3865658561145097
$ws.Range("P24").Value = 0.3865658561145097
$ws.Range("Q24").Value = 104.82920697909
$ws.Range("R24").Value = 943.4628628118098
$ws.Range("S24").Value = 0.001196481463182607
$ws.Range("T24").Value = 0.001196481463182607
$ws.Range("G25").Value = 0.628254
$ws.Range("H25").Value = 1.884762
$ws.Range("I25").Value = 0.003095155571184698
$ws.Range("J25").Value = 0.003095155571184698
$ws.Range("M25").Value = 41.09915599999999
$ws.Range("N25").Value = 123.297468
$ws.Range("O25").Value = 0.09521587377309249
$ws.Range("P25").Value = 0.09521587377309249
$ws.Range("Q25").Value = 25.82070915362399
$ws.Range("R25").Value = 232.3863823826159
$ws.Range("S25").Value = 0.0002947079421740061
$ws.Range("T25").Value = 0.0002947079421740062
$ws.Range("G26").Value = 0.628254
$ws.Range("H26").Value = 1.884762
$ws.Range("I26").Value = 0.003095155571184698
$ws.Range("J26").Value = 0.003095155571184698
$ws.Range("M26").Value = 112.0244103333333
$ws.Range("N26").Value = 336.073231
$ws.Range("O26").Value = 0.2595309284162377
$ws.Range("P26").Value = 0.2595309284162377
$ws.Range("Q26").Value = 70.379783889558
$ws.Range("R26").Value = 633.418055006022
$ws.Range("S26").Value = 0.0008032885989822551
$ws.Range("T26").Value = 0.0008032885989822552
